$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The four "previous state" tables on the diagram each contain a cell
# whose text begins with "prevTaskBook" (renamed to "prevOrganizer" in
# the underlying model, per the commit's doc/diagram rename pass).
$tableNames = @("Table 58", "Table 62", "Table 23", "Table 24")

foreach ($tableName in $tableNames) {
    $shp = $s.Shapes.Item($tableName)
    $tbl = $shp.Table

    for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
        for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
            $cellTextRange = $tbl.Cell($r, $c).Shape.TextFrame.TextRange

            if ($cellTextRange.Text -eq "prevTaskBook = s1") {
                # Single-paragraph cell ("prevTaskBook = s1") - assign the
                # leading run's new text directly so the trailing runs
                # (" " and "= s1") are left completely untouched.
                $cellTextRange.Text = "prevOrganizer"
            }
            else {
                $paragraphs = $cellTextRange.Paragraphs()
                for ($i = 1; $i -le $paragraphs.Count; $i++) {
                    $para = $paragraphs.Item($i)
                    if ($para.Text -eq "prevTaskBook = s3") {
                        $para.Text = "prevOrganizer = s3"
                    }
                }
            }
        }
    }
}
